$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new title row above the existing header row.
$ws.Rows("1:1").Insert()

# New title row content: A1/B1 hold the sheet title + version, C1 holds the source note.
$ws.Range("A1").Value = "NSIK <B>"
$ws.Range("B1").Value = "Versija: 1.0"
$ws.Range("C1").Value = "Šaltinis: LR AM įsakymas Nr. D1-346 (2024-10-28)"

# Formatting for the new title row: bold (inherited from the old header style),
# centered both ways for A1:B1, and wrapped text for the source note in C1.
$titleRng = $ws.Range("A1:B1")
$titleRng.Font.Bold = $true
$titleRng.HorizontalAlignment = -4108
$titleRng.VerticalAlignment = -4108

$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").WrapText = $true

$ws.Rows("1:1").RowHeight = 48

# Column widths: A/B tweak slightly, new column C is added.
$ws.Columns("A").ColumnWidth = 14.25
$ws.Columns("B").ColumnWidth = 47.083333333375776
$ws.Columns("C").ColumnWidth = 19.083333333401242

# Update the selection to match the new title row.
$ws.Range("A1:B1").Select()
